$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1305.5883
$ws.Range("I132").Value = 1324.6875
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3974.0625
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -1444.0625
$ws.Range("N132").Value = -8060

$ws.Range("H138").Value = 1933.7778
$ws.Range("I138").Value = 1073.8
$ws.Range("J138").Value = 3008.75
$ws.Range("K138").Value = 3221.4
$ws.Range("L138").Value = 9026.25
$ws.Range("M138").Value = 1918.6
$ws.Range("N138").Value = -19306.25

$ws.Range("H141").Value = 2337856.2
$ws.Range("I141").Value = 4670032
$ws.Range("J141").Value = 5680.8335
$ws.Range("K141").Value = 14010096
$ws.Range("L141").Value = 17042.5005
$ws.Range("M141").Value = -14004916
$ws.Range("N141").Value = -27402.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3323070.2
$ws.Range("I2").Value = 4651899
$ws.Range("J2").Value = 998.5
$ws.Range("K2").Value = 4651899
$ws.Range("L2").Value = 998.5
$ws.Range("M2").Value = -4651786
$ws.Range("N2").Value = -1224.5

$ws.Range("H32").Value = 3953.111
$ws.Range("I32").Value = 2242
$ws.Range("J32").Value = 11482
$ws.Range("K32").Value = 2242
$ws.Range("L32").Value = 11482
$ws.Range("M32").Value = -1955
$ws.Range("N32").Value = -12056

$ws.Range("H45").Value = 1719.2
$ws.Range("I45").Value = 1299.6666
$ws.Range("J45").Value = 1899
$ws.Range("K45").Value = 1299.6666
$ws.Range("L45").Value = 1899
$ws.Range("M45").Value = -922.6666
$ws.Range("N45").Value = -2653

$ws.Range("H61").Value = 2449.7173
$ws.Range("I61").Value = 1952.619
$ws.Range("J61").Value = 7669.25
$ws.Range("K61").Value = 1952.619
$ws.Range("L61").Value = 7669.25
$ws.Range("M61").Value = -1740.619
$ws.Range("N61").Value = -8093.25

$ws.Range("H97").Value = 693.7778
$ws.Range("I97").Value = 655.5
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 655.5
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -159.5
$ws.Range("N97").Value = -1992

$ws.Range("H116").Value = 3323070.2
$ws.Range("I116").Value = 4651899
$ws.Range("J116").Value = 998.5
$ws.Range("K116").Value = 4651899
$ws.Range("L116").Value = 998.5
$ws.Range("M116").Value = -4649605
$ws.Range("N116").Value = -5586.5

$ws.Range("H132").Value = 1356.5172
$ws.Range("I132").Value = 958.34784
$ws.Range("J132").Value = 2882.8333
$ws.Range("K132").Value = 2875.04352
$ws.Range("L132").Value = 8648.499899999999
$ws.Range("M132").Value = -345.0435200000002
$ws.Range("N132").Value = -13708.4999

$ws.Range("H136").Value = 2449.7173
$ws.Range("I136").Value = 1952.619
$ws.Range("J136").Value = 7669.25
$ws.Range("K136").Value = 5857.857
$ws.Range("L136").Value = 23007.75
$ws.Range("M136").Value = -3307.857
$ws.Range("N136").Value = -28107.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3323070.2
$ws.Range("I3").Value = 4651899
$ws.Range("J3").Value = 998.5
$ws.Range("K3").Value = 4651899
$ws.Range("L3").Value = 998.5
$ws.Range("M3").Value = -4651785
$ws.Range("N3").Value = -1226.5

$ws.Range("H86").Value = 93391.914
$ws.Range("I86").Value = 1649.5
$ws.Range("J86").Value = 221831.3
$ws.Range("K86").Value = 1649.5
$ws.Range("L86").Value = 221831.3
$ws.Range("M86").Value = -526.5
$ws.Range("N86").Value = -224077.3

$ws.Range("H89").Value = 93391.914
$ws.Range("I89").Value = 1649.5
$ws.Range("J89").Value = 221831.3
$ws.Range("K89").Value = 8247.5
$ws.Range("L89").Value = 1109156.5
$ws.Range("M89").Value = -2631.5
$ws.Range("N89").Value = -1120388.5

$ws.Range("H99").Value = 1424.5
$ws.Range("I99").Value = 999
$ws.Range("J99").Value = 1850
$ws.Range("K99").Value = 999
$ws.Range("L99").Value = 1850
$ws.Range("M99").Value = 499
$ws.Range("N99").Value = -4846

$ws.Range("H134").Value = 6906.6313
$ws.Range("I134").Value = 6757.7812
$ws.Range("J134").Value = 7700.5
$ws.Range("K134").Value = 20273.3436
$ws.Range("L134").Value = 23101.5
$ws.Range("M134").Value = -17738.3436
$ws.Range("N134").Value = -28171.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2630.3171
$ws.Range("I31").Value = 1691.5927
$ws.Range("J31").Value = 4440.7144
$ws.Range("K31").Value = 1691.5927
$ws.Range("L31").Value = 4440.7144
$ws.Range("M31").Value = -1396.5927
$ws.Range("N31").Value = -5030.7144

$ws.Range("H34").Value = 2630.3171
$ws.Range("I34").Value = 1691.5927
$ws.Range("J34").Value = 4440.7144
$ws.Range("K34").Value = 1691.5927
$ws.Range("L34").Value = 4440.7144
$ws.Range("M34").Value = -1489.5927
$ws.Range("N34").Value = -4844.7144

$ws.Range("H74").Value = 30000
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 30000
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31748

$ws.Range("H77").Value = 30000
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 30000
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -98736

$ws.Range("H107").Value = 431.7
$ws.Range("I107").Value = 401.8889
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 401.8889
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 1518.1111
$ws.Range("N107").Value = -4540

$ws.Range("H122").Value = 1914.3158
$ws.Range("I122").Value = 2025.6666
$ws.Range("J122").Value = 1723.4286
$ws.Range("K122").Value = 6076.9998
$ws.Range("L122").Value = 5170.2858
$ws.Range("M122").Value = -3626.9998
$ws.Range("N122").Value = -10070.2858

$ws.Range("H132").Value = 1887.0588
$ws.Range("I132").Value = 1098.2916
$ws.Range("J132").Value = 3780.1
$ws.Range("K132").Value = 3294.8748
$ws.Range("L132").Value = 11340.3
$ws.Range("M132").Value = -764.8748000000001
$ws.Range("N132").Value = -16400.3

$ws.Range("H134").Value = 1097
$ws.Range("I134").Value = 1097
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3291
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -756

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3500
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 4000
$ws.Range("K3").Value = 9000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = -8888
$ws.Range("N3").Value = -12224

$ws.Range("H26").Value = 320.8
$ws.Range("I26").Value = 299.5
$ws.Range("J26").Value = 335
$ws.Range("K26").Value = 898.5
$ws.Range("L26").Value = 1005
$ws.Range("M26").Value = -610.5
$ws.Range("N26").Value = -1581

$ws.Range("H107").Value = 1446.5883
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1446.5883
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4339.7649
$ws.Range("N107").Value = -8179.7649

$ws.Range("H113").Value = 7931.857
$ws.Range("I113").Value = 33685.668
$ws.Range("J113").Value = 908.0909
$ws.Range("K113").Value = 101057.004
$ws.Range("L113").Value = 2724.2727
$ws.Range("M113").Value = -98887.00399999999
$ws.Range("N113").Value = -7064.2727

$ws.Range("H129").Value = 91155.375
$ws.Range("I129").Value = 694
$ws.Range("J129").Value = 181616.75
$ws.Range("K129").Value = 2082
$ws.Range("L129").Value = 544850.25
$ws.Range("M129").Value = 2918
$ws.Range("N129").Value = -554850.25

$ws.Range("H131").Value = 17758.635
$ws.Range("I131").Value = 637.5
$ws.Range("J131").Value = 20693.686
$ws.Range("K131").Value = 1912.5
$ws.Range("L131").Value = 62081.058
$ws.Range("M131").Value = 3127.5
$ws.Range("N131").Value = -72161.058

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 15000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 15000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -15302

$ws.Range("H102").Value = 3237.923
$ws.Range("I102").Value = 4335.8
$ws.Range("J102").Value = 2551.75
$ws.Range("K102").Value = 4335.8
$ws.Range("L102").Value = 2551.75
$ws.Range("M102").Value = -2713.8
$ws.Range("N102").Value = -5795.75

$ws.Range("H122").Value = 1800
$ws.Range("I122").Value = 1800
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2950

$ws.Range("H126").Value = 2573597.8
$ws.Range("I126").Value = 2926850
$ws.Range("J126").Value = 336333.34
$ws.Range("K126").Value = 8780550
$ws.Range("L126").Value = 1009000.02
$ws.Range("M126").Value = -8778080
$ws.Range("N126").Value = -1013940.02

$ws.Range("H132").Value = 2407283.5
$ws.Range("I132").Value = 2749417.5
$ws.Range("J132").Value = 12344.5
$ws.Range("K132").Value = 8248252.5
$ws.Range("L132").Value = 37033.5
$ws.Range("M132").Value = -8245722.5
$ws.Range("N132").Value = -42093.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7100
$ws.Range("I40").Value = 1833.3334
$ws.Range("J40").Value = 15000
$ws.Range("K40").Value = 1833.3334
$ws.Range("L40").Value = 15000
$ws.Range("M40").Value = -1697.3334
$ws.Range("N40").Value = -15272

$ws.Range("H93").Value = 449.33334
$ws.Range("I93").Value = 320.57144
$ws.Range("J93").Value = 900
$ws.Range("K93").Value = 320.57144
$ws.Range("L93").Value = 900
$ws.Range("M93").Value = 927.4285600000001
$ws.Range("N93").Value = -3396

$ws.Range("H122").Value = 8961.799999999999
$ws.Range("I122").Value = 1804
$ws.Range("J122").Value = 10751.25
$ws.Range("K122").Value = 5412
$ws.Range("L122").Value = 32253.75
$ws.Range("M122").Value = -2962
$ws.Range("N122").Value = -37153.75

$ws.Range("H132").Value = 1808.0303
$ws.Range("I132").Value = 1588.1666
$ws.Range("J132").Value = 2071.8667
$ws.Range("K132").Value = 4764.4998
$ws.Range("L132").Value = 6215.6001
$ws.Range("M132").Value = -2234.4998
$ws.Range("N132").Value = -11275.6001

$ws.Range("H136").Value = 3279.2942
$ws.Range("I136").Value = 1749.7142
$ws.Range("J136").Value = 4350
$ws.Range("K136").Value = 5249.142599999999
$ws.Range("L136").Value = 13050
$ws.Range("M136").Value = -2699.142599999999
$ws.Range("N136").Value = -18150

$ws.Range("H140").Value = 54997.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54997.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54997.5
$ws.Range("N140").Value = -65357.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

$ws.Range("H80").Value = 79800
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 79800
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 79800
$ws.Range("N80").Value = -81796

$ws.Range("H83").Value = 79800
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 79800
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 239400
$ws.Range("N83").Value = -249384

$ws.Range("H126").Value = 8498.235000000001
$ws.Range("I126").Value = 8167
$ws.Range("J126").Value = 8870.875
$ws.Range("K126").Value = 24501
$ws.Range("L126").Value = 26612.625
$ws.Range("M126").Value = -22031
$ws.Range("N126").Value = -31552.625

$ws.Range("H135").Value = 85616
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 85616
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 85616
$ws.Range("N135").Value = -95756

$ws.Range("H136").Value = 14246935
$ws.Range("I136").Value = 19842928
$ws.Range("J136").Value = 2586.0908
$ws.Range("K136").Value = 59528784
$ws.Range("L136").Value = 7758.2724
$ws.Range("M136").Value = -59526234
$ws.Range("N136").Value = -12858.2724
